$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.917.57'
$ws.Range("E2").Value = '  -1.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.829.79'
$ws.Range("E3").Value = '  -1.53%  '

$ws.Range("E4").Value = '  +0.73%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.94'
$ws.Range("E5").Value = '  -1.00%  '

$ws.Range("E6").Value = '  +0.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4573'
$ws.Range("E7").Value = '  -0.95%  '

$ws.Range("E8").Value = '  -0.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07183'
$ws.Range("E9").Value = '  -2.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8766'
$ws.Range("E10").Value = '  -0.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07847'
$ws.Range("E11").Value = '  -0.07%  '

$ws.Range("E12").Value = '  -1.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.809.41'
$ws.Range("E13").Value = '  -2.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.337'
$ws.Range("E14").Value = '  -0.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.393'
$ws.Range("E15").Value = '  -2.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.10'
$ws.Range("E16").Value = '  -5.27%  '

$ws.Range("E17").Value = '  +0.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008728'
$ws.Range("E18").Value = '  -1.41%  '

$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.948.48'
$ws.Range("E20").Value = '  -1.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.50'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.003'
$ws.Range("E22").Value = '  -2.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.045.31'
$ws.Range("E23").Value = '  -5.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.43'
$ws.Range("E24").Value = '  -0.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.981'
$ws.Range("E25").Value = '  +4.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.23'
$ws.Range("E26").Value = '  -0.80%  '

$ws.Range("E27").Value = '  -0.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.967'
$ws.Range("E28").Value = '  -5.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.93'
$ws.Range("E29").Value = '  -1.83%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.931'
$ws.Range("E30").Value = '  -3.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08778'
$ws.Range("E31").Value = '  -0.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.055'
$ws.Range("E32").Value = '  +1.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7562'
$ws.Range("E33").Value = '  -0.61%  '

$ws.Range("E34").Value = '  -0.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.133'
$ws.Range("E35").Value = '  -3.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.574'
$ws.Range("E36").Value = '  -2.32%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.089'
$ws.Range("E37").Value = '  +1.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01937'
$ws.Range("E38").Value = '  -1.06%  '

$ws.Range("E39").Value = '  -1.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.901'
$ws.Range("E40").Value = '  -2.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.937'
$ws.Range("E41").Value = '  -1.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4978'
$ws.Range("E42").Value = '  -3.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1600'
$ws.Range("E43").Value = '  -2.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.309'
$ws.Range("E44").Value = '  -0.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4685'
$ws.Range("E45").Value = '  -3.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.008'
$ws.Range("E46").Value = '  +0.69%  '

$ws.Range("E47").Value = '  -1.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.23'
$ws.Range("E48").Value = '  -1.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.615'
$ws.Range("E49").Value = '  -2.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06117'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.45'
$ws.Range("E51").Value = '  -1.97%  '
